$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7 data: Username / Email / Joined.
# Shared strings are appended in first-use order, and the target workbook
# has "good@mail.com" (index 13) before "asdсавыфа" (index 14), so B7 must
# be written before A7.
$ws.Range("B7").Value = "good@mail.com"
$ws.Range("A7").Value = "asdсавыфа"
$ws.Range("C7").Value = 43136

# Turn B7 into a mailto: hyperlink (adds a new relationship + hyperlink entry)
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:good@mail.com")

# Hyperlinks.Add re-styles the cell with a freshly minted xf; put it back on
# the same "Гиперссылка" style the other email cells use (style index 2).
$ws.Range("B7").Style = $ws.Range("B6").Style

# Move the active selection to A7.
[void]$ws.Range("A7").Select()
